# "One more example completed"
#  - bump the cached "datetimeFigureOut" footer date from 11/02/2020 to
#    12/02/2020 everywhere it is cached (slide master, every slide layout,
#    and the notes master)
#  - flip the slide-1 title from "non-orderable" to "orderable"

$p = $ppt.ActivePresentation

$oldDate = "11/02/2020"
$newDate = "12/02/2020"

function Update-DateShape {
    param($container)

    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DateShape($p.SlideMaster)

# Every slide layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape($layouts.Item($li))
}

# Notes master footer date placeholder ("Espace réservé de la date 2").
# The notes master's placeholder shapes reject direct TextRange writes
# in this host, and routing the same update through
# NotesMaster.HeadersFooters.DateAndTime.Text instead lands on the
# wrong container after the CustomLayouts above have been walked (it
# clobbers slide-master shape 2), so it is intentionally left alone
# here rather than risk corrupting unrelated content.

# Slide 1 title: "non-orderable" -> "orderable".
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "Example E3: creatinine on 24h urine panel, non-orderable") {
            $tr.Text = "Example E3: creatinine on 24h urine panel, orderable"
        }
    }
}
